# Refactor, comment out the unnecessary
#
# - Rename the SKU codes from the old "QRTS0000N" scheme to the new
#   "KRTD0000N" scheme (shared-string values referenced by A2:A6).
# - Move the active selection to B4.
# - Re-fit the data columns (B:U, excluding the now default-width A/D)
#   to their new (slightly wider) column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename SKU codes -------------------------------------------------
$ws.Range("A2").Value = "KRTD00001"
$ws.Range("A3").Value = "KRTD00002"
$ws.Range("A4").Value = "KRTD00003"
$ws.Range("A5").Value = "KRTD00004"
$ws.Range("A6").Value = "KRTD00005"

# --- Re-fit column widths ----------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 31.084
$ws.Columns.Item(3).ColumnWidth = 24.584
$ws.Columns.Item(5).ColumnWidth = 15.417
$ws.Columns.Item(6).ColumnWidth = 34.084
$ws.Columns.Item(7).ColumnWidth = 19.25
$ws.Columns.Item(8).ColumnWidth = 18.417
$ws.Columns.Item(9).ColumnWidth = 19.584
$ws.Columns.Item(10).ColumnWidth = 18.417
$ws.Columns.Item(11).ColumnWidth = 14.084
$ws.Columns.Item(12).ColumnWidth = 20.25
$ws.Columns.Item(13).ColumnWidth = 21.75
$ws.Columns.Item(14).ColumnWidth = 17.917
$ws.Columns.Item(15).ColumnWidth = 17.917
$ws.Columns.Item(16).ColumnWidth = 14.084
$ws.Columns.Item(17).ColumnWidth = 22.25
$ws.Columns.Item(18).ColumnWidth = 21.584
$ws.Columns.Item(19).ColumnWidth = 21.084
$ws.Columns.Item(20).ColumnWidth = 19.417
$ws.Columns.Item(21).ColumnWidth = 18.25

# --- Move the active selection -----------------------------------------
$ws.Range("B4").Select()
